$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 195.40741
$ws.Range("I33").Value = 186.4
$ws.Range("J33").Value = 221.14285
$ws.Range("K33").Value = 186.4
$ws.Range("L33").Value = 221.14285
$ws.Range("M33").Value = 42.59999999999999
$ws.Range("N33").Value = -679.14285
$ws.Range("H74").Value = 6969.4062
$ws.Range("J74").Value = 7067.3667
$ws.Range("L74").Value = 7067.3667
$ws.Range("N74").Value = -8939.366699999999
$ws.Range("H77").Value = 6969.4062
$ws.Range("J77").Value = 7067.3667
$ws.Range("L77").Value = 35336.8335
$ws.Range("N77").Value = -44696.8335
$ws.Range("H92").Value = 819.9677
$ws.Range("I92").Value = 263.81818
$ws.Range("J92").Value = 2179.4443
$ws.Range("K92").Value = 263.81818
$ws.Range("L92").Value = 2179.4443
$ws.Range("M92").Value = 984.18182
$ws.Range("N92").Value = -4675.4443
$ws.Range("H99").Value = 1251.7646
$ws.Range("I99").Value = 350.4
$ws.Range("J99").Value = 1627.3334
$ws.Range("K99").Value = 1051.2
$ws.Range("L99").Value = 4882.0002
$ws.Range("M99").Value = 446.8000000000002
$ws.Range("N99").Value = -7878.0002
$ws.Range("H101").Value = 1165.8889
$ws.Range("I101").Value = 732.7143
$ws.Range("J101").Value = 2682
$ws.Range("K101").Value = 2198.1429
$ws.Range("L101").Value = 8046
$ws.Range("M101").Value = -576.1428999999998
$ws.Range("N101").Value = -11290
$ws.Range("H116").Value = 55298.69
$ws.Range("I116").Value = 110833.336
$ws.Range("J116").Value = 7697.5713
$ws.Range("K116").Value = 110833.336
$ws.Range("L116").Value = 7697.5713
$ws.Range("M116").Value = -107391.336
$ws.Range("N116").Value = -14581.5713
$ws.Range("H127").Value = 964.6667
$ws.Range("I127").Value = 964.6667
$ws.Range("K127").Value = 2894.0001
$ws.Range("M127").Value = 2065.9999
$ws.Range("H137").Value = 3299.7856
$ws.Range("I137").Value = 1589.5
$ws.Range("J137").Value = 3584.8333
$ws.Range("K137").Value = 4768.5
$ws.Range("L137").Value = 10754.4999
$ws.Range("M137").Value = -2218.5
$ws.Range("N137").Value = -15854.4999
$ws.Range("H138").Value = 2451.8262
$ws.Range("J138").Value = 2870.3667
$ws.Range("L138").Value = 8611.1001
$ws.Range("N138").Value = -18891.1001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2026.641
$ws.Range("I74").Value = 1480.3125
$ws.Range("J74").Value = 4524.143
$ws.Range("K74").Value = 1480.3125
$ws.Range("L74").Value = 4524.143
$ws.Range("M74").Value = -606.3125
$ws.Range("N74").Value = -6272.143
$ws.Range("H77").Value = 2026.641
$ws.Range("I77").Value = 1480.3125
$ws.Range("J77").Value = 4524.143
$ws.Range("K77").Value = 7401.5625
$ws.Range("L77").Value = 22620.715
$ws.Range("M77").Value = -3033.5625
$ws.Range("N77").Value = -31356.715
$ws.Range("H102").Value = 3405.3
$ws.Range("I102").Value = 3231.6875
$ws.Range("J102").Value = 4099.75
$ws.Range("K102").Value = 3231.6875
$ws.Range("L102").Value = 4099.75
$ws.Range("M102").Value = -1609.6875
$ws.Range("N102").Value = -7343.75
$ws.Range("H122").Value = 2603.75
$ws.Range("I122").Value = 2401.1155
$ws.Range("K122").Value = 7203.3465
$ws.Range("M122").Value = -4753.3465
$ws.Range("H132").Value = 6537.5293
$ws.Range("I132").Value = 4595
$ws.Range("K132").Value = 13785
$ws.Range("M132").Value = -11255

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2530.4
$ws.Range("I20").Value = 2782
$ws.Range("J20").Value = 2362.6667
$ws.Range("K20").Value = 2782
$ws.Range("L20").Value = 2362.6667
$ws.Range("M20").Value = -2535
$ws.Range("N20").Value = -2856.6667
$ws.Range("H107").Value = 2594.04
$ws.Range("I107").Value = 1948.5555
$ws.Range("K107").Value = 1948.5555
$ws.Range("M107").Value = -28.55549999999994
$ws.Range("H133").Value = 85601.5
$ws.Range("J133").Value = 99997
$ws.Range("L133").Value = 99997
$ws.Range("N133").Value = -110117

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4817.096
$ws.Range("I99").Value = 4706.294
$ws.Range("J99").Value = 5026.3887
$ws.Range("K99").Value = 4706.294
$ws.Range("L99").Value = 5026.3887
$ws.Range("M99").Value = -3208.294
$ws.Range("N99").Value = -8022.3887
$ws.Range("H122").Value = 2061.652
$ws.Range("I122").Value = 1595.9
$ws.Range("J122").Value = 5166.6665
$ws.Range("K122").Value = 4787.700000000001
$ws.Range("L122").Value = 15499.9995
$ws.Range("M122").Value = -2337.700000000001
$ws.Range("N122").Value = -20399.9995
$ws.Range("H126").Value = 4817.096
$ws.Range("I126").Value = 4706.294
$ws.Range("J126").Value = 5026.3887
$ws.Range("K126").Value = 14118.882
$ws.Range("L126").Value = 15079.1661
$ws.Range("M126").Value = -11648.882
$ws.Range("N126").Value = -20019.1661

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1296.25
$ws.Range("J113").Value = 1296.25
$ws.Range("L113").Value = 3888.75
$ws.Range("N113").Value = -8228.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 16658286
$ws.Range("I11").Value = 16658286
$ws.Range("K11").Value = 16658286
$ws.Range("M11").Value = -16658147
$ws.Range("H31").Value = 5928.75
$ws.Range("I31").Value = 3918.5715
$ws.Range("K31").Value = 3918.5715
$ws.Range("M31").Value = -3626.5715
$ws.Range("H37").Value = 5928.75
$ws.Range("I37").Value = 3918.5715
$ws.Range("K37").Value = 3918.5715
$ws.Range("M37").Value = -3641.5715
$ws.Range("H136").Value = 26290.8
$ws.Range("J136").Value = 26290.8
$ws.Range("L136").Value = 78872.39999999999
$ws.Range("N136").Value = -83972.39999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1866.6666
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5340
$ws.Range("H22").Value = 101173.6
$ws.Range("I22").Value = 143676.58
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 143676.58
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -143381.58
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 101173.6
$ws.Range("I27").Value = 143676.58
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 143676.58
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -143569.58
$ws.Range("N27").Value = -2214
$ws.Range("H43").Value = 503000
$ws.Range("J43").Value = 1000000
$ws.Range("L43").Value = 1000000
$ws.Range("N43").Value = -1000386
$ws.Range("H46").Value = 1144.4286
$ws.Range("I46").Value = 1144.4286
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1144.4286
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -956.4286
$ws.Range("H68").Value = 4125.55
$ws.Range("I68").Value = 2529.8
$ws.Range("J68").Value = 5721.3
$ws.Range("K68").Value = 2529.8
$ws.Range("L68").Value = 5721.3
$ws.Range("M68").Value = -1780.8
$ws.Range("N68").Value = -7219.3
$ws.Range("H71").Value = 4125.55
$ws.Range("I71").Value = 2529.8
$ws.Range("J71").Value = 5721.3
$ws.Range("K71").Value = 12649
$ws.Range("L71").Value = 28606.5
$ws.Range("M71").Value = -8905
$ws.Range("N71").Value = -36094.5
$ws.Range("H136").Value = 8295.727999999999
$ws.Range("I136").Value = 7531.6875
$ws.Range("K136").Value = 22595.0625
$ws.Range("M136").Value = -20045.0625
$ws.Range("N46").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("H58").Value = 2060
$ws.Range("I58").Value = 2060
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2060
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1752
$ws.Range("H132").Value = 4183.6743
$ws.Range("J132").Value = 6180.6
$ws.Range("L132").Value = 18541.8
$ws.Range("N132").Value = -23601.8
$ws.Range("H136").Value = 6908958
$ws.Range("I136").Value = 8224488.5
$ws.Range("K136").Value = 24673465.5
$ws.Range("M136").Value = -24670915.5
$ws.Range("M44").ClearContents()
$ws.Range("N58").ClearContents()
